# Regenerate merged AHB files
# - rename the "*_old" header columns to "*_FV2410"
# - rename the "*_new" header columns to "*_FV2504"
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (row 1) shared strings -----------------------
# Columns A..J were "<Name>_old" -> "<Name>_FV2410"
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2410Headers[$i]
}

# Column K ("diff") is untouched.

# Columns L..U were "<Name>_new" -> "<Name>_FV2504"
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)
for ($i = 0; $i -lt $fv2504Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2504Headers[$i]
}

# --- 2) Turn the whole data range into a real Excel Table -------------------
$dataRange = $ws.UsedRange
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3) Freeze the header row ------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Applied AHB merge regeneration: renamed headers, added Table1, froze header row."
